$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the counts (B=in_deck, C=in_reserve, D=in_discard) for several cities
$ws.Range("B13").Value = 0

$ws.Range("B18").Value = 2
$ws.Range("C18").Value = 0

$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 1

$ws.Range("B30").Value = 2

# Reset the view: scroll back to top-left and select A1
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
